$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.36022366666667
$ws.Range("H2").Value = 58.080671
$ws.Range("I2").Value = 0.005884129141485179
$ws.Range("J2").Value = 0.005884129141485179
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 73.15142415085633
$ws.Range("R2").Value = 658.3628173577069
$ws.Range("S2").Value = 0.002502503370772032
$ws.Range("T2").Value = 0.002502503370772032
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.36022366666667
$ws.Range("H3").Value = 58.080671
$ws.Range("I3").Value = 0.005884129141485179
$ws.Range("J3").Value = 0.005884129141485179
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 83.89565777121221
$ws.Range("R3").Value = 755.0609199409099
$ws.Range("S3").Value = 0.002870062596903489
$ws.Range("T3").Value = 0.002870062596903489
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.36022366666667
$ws.Range("H4").Value = 58.080671
$ws.Range("I4").Value = 0.005884129141485179
$ws.Range("J4").Value = 0.005884129141485179
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 5.633541137052889
$ws.Range("R4").Value = 50.701870233476
$ws.Range("S4").Value = 0.0001927229148100287
$ws.Range("T4").Value = 0.0001927229148100287
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.36022366666667
$ws.Range("H5").Value = 58.080671
$ws.Range("I5").Value = 0.005884129141485179
$ws.Range("J5").Value = 0.005884129141485179
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 9.320114927659555
$ws.Range("R5").Value = 83.88103434893598
$ws.Range("S5").Value = 0.0003188402589996292
$ws.Range("T5").Value = 0.0003188402589996292
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3161.845459
$ws.Range("H6").Value = 9485.536377
$ws.Range("I6").Value = 0.9609758299542277
$ws.Range("J6").Value = 0.9609758299542278
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 11946.8401942585
$ws.Range("R6").Value = 107521.5617483265
$ws.Range("S6").Value = 0.4087002844203233
$ws.Range("T6").Value = 0.4087002844203234
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3161.845459
$ws.Range("H7").Value = 9485.536377
$ws.Range("I7").Value = 0.9609758299542277
$ws.Range("J7").Value = 0.9609758299542278
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 13701.55165151546
$ws.Range("R7").Value = 123313.9648636392
$ws.Range("S7").Value = 0.4687287990731914
$ws.Range("T7").Value = 0.4687287990731915
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3161.845459
$ws.Range("H8").Value = 9485.536377
$ws.Range("I8").Value = 0.9609758299542277
$ws.Range("J8").Value = 0.9609758299542278
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 920.0506548356014
$ws.Range("R8").Value = 8280.455893520413
$ws.Range("S8").Value = 0.03147484675430143
$ws.Range("T8").Value = 0.03147484675430143
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3161.845459
$ws.Range("H9").Value = 9485.536377
$ws.Range("I9").Value = 0.9609758299542277
$ws.Range("J9").Value = 0.9609758299542278
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 1522.129267138381
$ws.Range("R9").Value = 13699.16340424543
$ws.Range("S9").Value = 0.05207189970641153
$ws.Range("T9").Value = 0.05207189970641154
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.055785333333333
$ws.Range("H10").Value = 6.167356
$ws.Range("I10").Value = 0.0006248123263850286
$ws.Range("J10").Value = 0.0006248123263850286
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 7.767659479094666
$ws.Range("R10").Value = 69.90893531185199
$ws.Range("S10").Value = 0.0002657309034661654
$ws.Range("T10").Value = 0.0002657309034661654
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.055785333333333
$ws.Range("H11").Value = 6.167356
$ws.Range("I11").Value = 0.0006248123263850286
$ws.Range("J11").Value = 0.0006248123263850286
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 8.908547016084443
$ws.Range("R11").Value = 80.17692314476
$ws.Range("S11").Value = 0.0003047605592812162
$ws.Range("T11").Value = 0.0003047605592812162
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.055785333333333
$ws.Range("H12").Value = 6.167356
$ws.Range("I12").Value = 0.0006248123263850286
$ws.Range("J12").Value = 0.0006248123263850286
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 0.5982033804817778
$ws.Range("R12").Value = 5.383830424336001
$ws.Range("S12").Value = 0.0000204644816343654
$ws.Range("T12").Value = 0.0000204644816343654
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.055785333333333
$ws.Range("H13").Value = 6.167356
$ws.Range("I13").Value = 0.0006248123263850286
$ws.Range("J13").Value = 0.0006248123263850286
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 0.9896660236551109
$ws.Range("R13").Value = 8.906994212895999
$ws.Range("S13").Value = 0.00003385638200328156
$ws.Range("T13").Value = 0.00003385638200328156
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 106.9830526666667
$ws.Range("H14").Value = 320.949158
$ws.Range("I14").Value = 0.03251522857790212
$ws.Range("J14").Value = 0.03251522857790212
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.778439
$ws.Range("N14").Value = 11.335317
$ws.Range("O14").Value = 0.4252971528324392
$ws.Range("P14").Value = 0.4252971528324392
$ws.Range("Q14").Value = 404.2289385347873
$ws.Range("R14").Value = 3638.060446813086
$ws.Range("S14").Value = 0.01382863413787773
$ws.Range("T14").Value = 0.01382863413787773
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 106.9830526666667
$ws.Range("H15").Value = 320.949158
$ws.Range("I15").Value = 0.03251522857790212
$ws.Range("J15").Value = 0.03251522857790212
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.4877633593505858
$ws.Range("P15").Value = 0.4877633593505858
$ws.Range("Q15").Value = 463.6007170359089
$ws.Range("R15").Value = 4172.406453323179
$ws.Range("S15").Value = 0.01585973712120971
$ws.Range("T15").Value = 0.01585973712120971
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 106.9830526666667
$ws.Range("H16").Value = 320.949158
$ws.Range("I16").Value = 0.03251522857790212
$ws.Range("J16").Value = 0.03251522857790212
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.2909853333333334
$ws.Range("N16").Value = 0.8729560000000001
$ws.Range("O16").Value = 0.03275300561492853
$ws.Range("P16").Value = 0.03275300561492853
$ws.Range("Q16").Value = 31.13049924122756
$ws.Range("R16").Value = 280.174493171048
$ws.Range("S16").Value = 0.001064971464182713
$ws.Range("T16").Value = 0.001064971464182713
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 106.9830526666667
$ws.Range("H17").Value = 320.949158
$ws.Range("I17").Value = 0.03251522857790212
$ws.Range("J17").Value = 0.03251522857790212
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.4814053333333333
$ws.Range("N17").Value = 1.444216
$ws.Range("O17").Value = 0.0541864822020464
$ws.Range("P17").Value = 0.05418648220204641
$ws.Range("Q17").Value = 51.50221213001421
$ws.Range("R17").Value = 463.519909170128
$ws.Range("S17").Value = 0.001761885854631964
$ws.Range("T17").Value = 0.001761885854631964
